# Kansas overview workbook - convert numeric "count" cells to text cells
# (same displayed value) and refresh the all-zero County rows to the new
# percent/currency text formatting, plus append a "Total" row to County.
#
# Technique: write a `=TEXT(<value>,"0")` (or literal `="..."` ) formula
# into the cell, then Copy + PasteSpecial(xlPasteValues = -4163) over
# itself. That converts the formula result into a literal stored value
# without touching NumberFormat/style (so no new cellXfs entries are
# created) and yields a genuine text-typed cell, matching the original
# file's inline-string cells.

$wb = $excel.ActiveWorkbook

$xlPasteValues = -4163

function Convert-NumberCellToText($cell) {
    $v = $cell.Value2
    $cell.Formula = "=TEXT($v,""0"")"
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

function Set-LiteralText($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = "=""$escaped"""
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

# ---------------------------------------------------------------------
# Sheet "Overall": A2 (No. of 990 Filers w/ Gov Grants) -> text "971"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Convert-NumberCellToText $wsOverall.Cells.Item(2, 1)

# ---------------------------------------------------------------------
# Sheet "County": column B rows 2-89 -> text (same numeric value)
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
for ($r = 2; $r -le 89; $r++) {
    Convert-NumberCellToText $wsCounty.Cells.Item($r, 2)
}

# Rows 90-98 are the all-zero counties; refresh to formatted text values.
for ($r = 90; $r -le 98; $r++) {
    Set-LiteralText $wsCounty.Cells.Item($r, 2) "0.00%"
    Set-LiteralText $wsCounty.Cells.Item($r, 3) "`$0"
    Set-LiteralText $wsCounty.Cells.Item($r, 4) "0.00%"
    Set-LiteralText $wsCounty.Cells.Item($r, 5) "0.00%"
    Set-LiteralText $wsCounty.Cells.Item($r, 6) "0.00%"
}

# Append the new "Total" row (row 99).
Set-LiteralText $wsCounty.Cells.Item(99, 1) "Total"
Set-LiteralText $wsCounty.Cells.Item(99, 2) "971"
Set-LiteralText $wsCounty.Cells.Item(99, 3) "`$1,044,083,084"
Set-LiteralText $wsCounty.Cells.Item(99, 4) "8.02%"
Set-LiteralText $wsCounty.Cells.Item(99, 5) "-11.92%"
Set-LiteralText $wsCounty.Cells.Item(99, 6) "70.13%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": column B rows 2-6 -> text
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
for ($r = 2; $r -le 6; $r++) {
    Convert-NumberCellToText $wsCd.Cells.Item($r, 2)
}

# ---------------------------------------------------------------------
# Sheet "Size": column B rows 2-8 -> text
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
for ($r = 2; $r -le 8; $r++) {
    Convert-NumberCellToText $wsSize.Cells.Item($r, 2)
}

# ---------------------------------------------------------------------
# Sheet "Subsector": column B rows 2-13 -> text
# ---------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")
for ($r = 2; $r -le 13; $r++) {
    Convert-NumberCellToText $wsSubsector.Cells.Item($r, 2)
}
